$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (id 111618070 -> 111618078, takes values from old row 5)
$ws.Range("A2").Value2 = 111618078
$ws.Range("I2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("Q2").Value2 = 580612.1009209087
$ws.Range("R2").Value2 = 6415119.491031807
$ws.Range("AC2").ClearContents()

# Row 4 (id 111618046 -> 111618039, takes values from old row 9)
$ws.Range("A4").Value2 = 111618039
$ws.Range("Q4").Value2 = 580599.6803078586
$ws.Range("R4").Value2 = 6415233.627682217

# Row 5 (id 111618078 -> 111618070, takes values from old row 2)
$ws.Range("A5").Value2 = 111618070
$ws.Range("I5").Value2 = "15"
$ws.Range("K5").Value2 = "blomning"
$ws.Range("Q5").Value2 = 580592.470229132
$ws.Range("R5").Value2 = 6415141.442167919
$ws.Range("AC5").Value2 = "1 blomma"

# Row 6 (id 111618056 -> 111618109, takes values from old row 8)
$ws.Range("A6").Value2 = 111618109
$ws.Range("I6").Value2 = "10"
$ws.Range("P6").Value2 = "A 32649, Sm"
$ws.Range("Q6").Value2 = 580619.1666838422
$ws.Range("R6").Value2 = 6415112.716507593
$ws.Range("AC6").Value2 = "1 blomma"

# Row 7 (id 111618144 -> 111618056, takes values from old row 6)
$ws.Range("A7").Value2 = 111618056
$ws.Range("I7").Value2 = "15"
$ws.Range("K7").Value2 = "blomning"
$ws.Range("Q7").Value2 = 580582.6881743574
$ws.Range("R7").Value2 = 6415124.22061418
$ws.Range("AC7").Value2 = "2 blommor"

# Row 8 (id 111618109 -> 111618144, takes values from old row 7)
$ws.Range("A8").Value2 = 111618144
$ws.Range("I8").Value2 = "2"
$ws.Range("K8").ClearContents()
$ws.Range("P8").Value2 = "A 32649, Heda, Sm"
$ws.Range("Q8").Value2 = 580620.6996611424
$ws.Range("R8").Value2 = 6415142.541277731
$ws.Range("AC8").ClearContents()

# Row 9 (id 111618039 -> 111618046, takes values from old row 4)
$ws.Range("A9").Value2 = 111618046
$ws.Range("Q9").Value2 = 580591.6383206119
$ws.Range("R9").Value2 = 6415156.322361182

Write-Host "done"
